$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 201.25
$ws.Range("I11").Value = 201.25
$ws.Range("K11").Value = 201.25
$ws.Range("M11").Value = -61.25
$ws.Range("H17").Value = 680.7857
$ws.Range("I17").Value = 392.2
$ws.Range("J17").Value = 841.1111
$ws.Range("K17").Value = 1176.6
$ws.Range("L17").Value = 2523.3333
$ws.Range("M17").Value = -1008.6
$ws.Range("N17").Value = -2859.3333
$ws.Range("H33").Value = 306.6
$ws.Range("I33").Value = 266.57144
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 266.57144
$ws.Range("L33").Value = 400
$ws.Range("M33").Value = -37.57144
$ws.Range("N33").Value = -858
$ws.Range("H43").Value = 6082.7144
$ws.Range("I43").Value = 3616
$ws.Range("K43").Value = 3616
$ws.Range("M43").Value = -3547
$ws.Range("H74").Value = 3000000
$ws.Range("I74").Value = 3000000
$ws.Range("K74").Value = 3000000
$ws.Range("M74").Value = -2999064
$ws.Range("H77").Value = 3000000
$ws.Range("I77").Value = 3000000
$ws.Range("K77").Value = 15000000
$ws.Range("M77").Value = -14995320
$ws.Range("H99").Value = 500
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 500
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1500
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -4496
$ws.Range("H106").Value = 90946456
$ws.Range("I106").Value = 100033600
$ws.Range("K106").Value = 100033600
$ws.Range("M106").Value = -100032969
$ws.Range("H112").Value = 1063.9
$ws.Range("J112").Value = 1054.3334
$ws.Range("L112").Value = 3163.0002
$ws.Range("N112").Value = -5379.0002
$ws.Range("H132").Value = 977.8
$ws.Range("I132").Value = 977.8
$ws.Range("K132").Value = 2933.4
$ws.Range("M132").Value = -403.3999999999996
$ws.Range("H138").Value = 2211.5625
$ws.Range("I138").Value = 1661.5
$ws.Range("K138").Value = 4984.5
$ws.Range("M138").Value = 155.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5551.8135
$ws.Range("I32").Value = 2762.3635
$ws.Range("J32").Value = 13734.2
$ws.Range("K32").Value = 2762.3635
$ws.Range("L32").Value = 13734.2
$ws.Range("M32").Value = -2475.3635
$ws.Range("N32").Value = -14308.2
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -788
$ws.Range("H88").Value = 825.7857
$ws.Range("I88").Value = 561.5
$ws.Range("J88").Value = 931.5
$ws.Range("K88").Value = 561.5
$ws.Range("L88").Value = 931.5
$ws.Range("M88").Value = -155.5
$ws.Range("N88").Value = -1743.5
$ws.Range("H91").Value = 825.7857
$ws.Range("I91").Value = 561.5
$ws.Range("J91").Value = 931.5
$ws.Range("K91").Value = 561.5
$ws.Range("L91").Value = 931.5
$ws.Range("M91").Value = 842.5
$ws.Range("N91").Value = -3739.5
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 19598.834
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 19598.834
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 19598.834
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -20218.834
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("H105").Value = 3626366.8
$ws.Range("J105").Value = 3499.75
$ws.Range("L105").Value = 3499.75
$ws.Range("N105").Value = -6993.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5580
$ws.Range("I31").Value = 4492.75
$ws.Range("J31").Value = 6123.625
$ws.Range("K31").Value = 4492.75
$ws.Range("L31").Value = 6123.625
$ws.Range("M31").Value = -4197.75
$ws.Range("N31").Value = -6713.625
$ws.Range("H34").Value = 5580
$ws.Range("I34").Value = 4492.75
$ws.Range("J34").Value = 6123.625
$ws.Range("K34").Value = 4492.75
$ws.Range("L34").Value = 6123.625
$ws.Range("M34").Value = -4290.75
$ws.Range("N34").Value = -6527.625
$ws.Range("H68").Value = 49999
$ws.Range("J68").Value = 49999
$ws.Range("L68").Value = 49999
$ws.Range("N68").Value = -51497
$ws.Range("H71").Value = 49999
$ws.Range("J71").Value = 49999
$ws.Range("L71").Value = 149997
$ws.Range("N71").Value = -157485
$ws.Range("H99").Value = 12150.593
$ws.Range("I99").Value = 8647.571
$ws.Range("J99").Value = 15923.077
$ws.Range("K99").Value = 8647.571
$ws.Range("L99").Value = 15923.077
$ws.Range("M99").Value = -7149.571
$ws.Range("N99").Value = -18919.077
$ws.Range("H107").Value = 31250402
$ws.Range("I107").Value = 33333746
$ws.Range("K107").Value = 33333746
$ws.Range("M107").Value = -33331826
$ws.Range("H122").Value = 924.55554
$ws.Range("I122").Value = 960.2857
$ws.Range("J122").Value = 799.5
$ws.Range("K122").Value = 2880.8571
$ws.Range("L122").Value = 2398.5
$ws.Range("M122").Value = -430.8571000000002
$ws.Range("N122").Value = -7298.5
$ws.Range("H126").Value = 12150.593
$ws.Range("I126").Value = 8647.571
$ws.Range("J126").Value = 15923.077
$ws.Range("K126").Value = 25942.713
$ws.Range("L126").Value = 47769.231
$ws.Range("M126").Value = -23472.713
$ws.Range("N126").Value = -52709.231
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 592.1622
$ws.Range("I4").Value = 545.96875
$ws.Range("K4").Value = 1637.90625
$ws.Range("M4").Value = -1525.90625
$ws.Range("H23").Value = 166721.67
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 200062
$ws.Range("K23").Value = 60
$ws.Range("L23").Value = 600186
$ws.Range("M23").Value = 175
$ws.Range("N23").Value = -600656
$ws.Range("H62").Value = 2125
$ws.Range("J62").Value = 2125
$ws.Range("L62").Value = 6375
$ws.Range("N62").Value = -7747
$ws.Range("H65").Value = 2125
$ws.Range("J65").Value = 2125
$ws.Range("L65").Value = 19125
$ws.Range("N65").Value = -25989
$ws.Range("H97").Value = 4726.3335
$ws.Range("I97").Value = 4594.5
$ws.Range("J97").Value = 4990
$ws.Range("K97").Value = 13783.5
$ws.Range("L97").Value = 14970
$ws.Range("M97").Value = -13287.5
$ws.Range("N97").Value = -15962
$ws.Range("H107").Value = 63056.188
$ws.Range("I107").Value = 461.25
$ws.Range("J107").Value = 83921.164
$ws.Range("K107").Value = 1383.75
$ws.Range("L107").Value = 251763.492
$ws.Range("M107").Value = 536.25
$ws.Range("N107").Value = -255603.492
$ws.Range("H109").Value = 994.1177
$ws.Range("I109").Value = 994.1177
$ws.Range("K109").Value = 2982.3531
$ws.Range("M109").Value = -1942.3531
$ws.Range("H121").Value = 489.625
$ws.Range("I121").Value = 348.5
$ws.Range("J121").Value = 630.75
$ws.Range("K121").Value = 1045.5
$ws.Range("L121").Value = 1892.25
$ws.Range("M121").Value = 264.5
$ws.Range("N121").Value = -4512.25
$ws.Range("H131").Value = 1728.125
$ws.Range("J131").Value = 1728.125
$ws.Range("L131").Value = 5184.375
$ws.Range("N131").Value = -15264.375
$ws.Range("H137").Value = 5147.143
$ws.Range("J137").Value = 5000
$ws.Range("L137").Value = 15000
$ws.Range("N137").Value = -25200
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 35751710
$ws.Range("I113").Value = 62520500
$ws.Range("K113").Value = 62520500
$ws.Range("M113").Value = -62518330
$ws.Range("H122").Value = 36639.516
$ws.Range("I122").Value = 1937.45
$ws.Range("K122").Value = 5812.35
$ws.Range("M122").Value = -3362.35
$ws.Range("H132").Value = 1802.9474
$ws.Range("I132").Value = 1537.7333
$ws.Range("J132").Value = 2797.5
$ws.Range("K132").Value = 4613.199900000001
$ws.Range("L132").Value = 8392.5
$ws.Range("M132").Value = -2083.199900000001
$ws.Range("N132").Value = -13452.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2705
$ws.Range("I31").Value = 1426.6666
$ws.Range("K31").Value = 1426.6666
$ws.Range("M31").Value = -1178.6666
$ws.Range("H43").Value = 13346
$ws.Range("I43").Value = 9999
$ws.Range("J43").Value = 15019.5
$ws.Range("K43").Value = 9999
$ws.Range("L43").Value = 15019.5
$ws.Range("M43").Value = -9806
$ws.Range("N43").Value = -15405.5
$ws.Range("H46").Value = 86749.5
$ws.Range("I46").Value = 2999.2856
$ws.Range("J46").Value = 203999.8
$ws.Range("K46").Value = 2999.2856
$ws.Range("L46").Value = 203999.8
$ws.Range("M46").Value = -2811.2856
$ws.Range("N46").Value = -204375.8
$ws.Range("H61").Value = 12347757
$ws.Range("I61").Value = 18520168
$ws.Range("J61").Value = 2932.3333
$ws.Range("K61").Value = 18520168
$ws.Range("L61").Value = 2932.3333
$ws.Range("M61").Value = -18519966
$ws.Range("N61").Value = -3336.3333
$ws.Range("H113").Value = 12347757
$ws.Range("I113").Value = 18520168
$ws.Range("J113").Value = 2932.3333
$ws.Range("K113").Value = 18520168
$ws.Range("L113").Value = 2932.3333
$ws.Range("M113").Value = -18517998
$ws.Range("N113").Value = -7272.3333
$ws.Range("H122").Value = 2806.5715
$ws.Range("I122").Value = 2441
$ws.Range("K122").Value = 7323
$ws.Range("M122").Value = -4873
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 29999.5
$ws.Range("J104").Value = 29999.5
$ws.Range("L104").Value = 29999.5
$ws.Range("N104").Value = -36987.5
$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960
$ws.Range("H132").Value = 1732.9286
$ws.Range("I132").Value = 1846.5
$ws.Range("J132").Value = 1051.5
$ws.Range("K132").Value = 5539.5
$ws.Range("L132").Value = 3154.5
$ws.Range("M132").Value = -3009.5
$ws.Range("N132").Value = -8214.5
